$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.225.09"
$ws.Range("E2").Value = "  -0.63%  "
$ws.Range("D3").Value = "1.839.08"
$ws.Range("E3").Value = "  -1.37%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "232.70"
$ws.Range("E5").Value = "  -1.16%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4657"
$ws.Range("E7").Value = "  -3.09%  "
$ws.Range("E8").Value = "  -2.46%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06269"
$ws.Range("E9").Value = "  -3.96%  "
$ws.Range("D10").Value = "1.830.50"
$ws.Range("E10").Value = "  -2.00%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07421"
$ws.Range("E11").Value = "  -0.27%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "16.10"
$ws.Range("E12").Value = "  -0.38%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.933"
$ws.Range("E13").Value = "  -2.76%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "83.50"
$ws.Range("E14").Value = "  -3.91%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6186"
$ws.Range("E15").Value = "  -3.00%  "
$ws.Range("D16").Value = "30.159.17"
$ws.Range("E16").Value = "  -0.78%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.000"
$ws.Range("E17").Value = "  +0.01%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "225.69"
$ws.Range("E18").Value = "  -2.72%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000007270"
$ws.Range("E19").Value = "  -2.53%  "
$ws.Range("B20").Value = "BinanceUSD"
$ws.Range("C20").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "1.001"
$ws.Range("E20").Value = "  -0.02%  "
$ws.Range("B21").Value = "Avalanche"
$ws.Range("C21").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.30"
$ws.Range("E21").Value = "  -5.07%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.877"
$ws.Range("E22").Value = "  -4.81%  "
$ws.Range("E23").Value = "  -3.90%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "9.184"
$ws.Range("E24").Value = "  -1.29%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "164.17"
$ws.Range("E25").Value = "  -2.85%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "17.76"
$ws.Range("E26").Value = "  -2.01%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.859"
$ws.Range("E27").Value = "  -1.84%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.1031"
$ws.Range("E28").Value = "  -1.35%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.373"
$ws.Range("E29").Value = "  -0.46%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.064"
$ws.Range("E30").Value = "  -4.40%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.796"
$ws.Range("E31").Value = "  -4.24%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.04817"
$ws.Range("E32").Value = "  -2.96%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.138"
$ws.Range("E33").Value = "  -2.43%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7023"
$ws.Range("E34").Value = "  -4.94%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.698"
$ws.Range("E35").Value = "  -0.49%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.01870"
$ws.Range("E36").Value = "  -3.37%  "
$ws.Range("E37").Value = "  +0.69%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.8895"
$ws.Range("E38").Value = "  -2.55%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "104.50"
$ws.Range("E39").Value = "  -1.53%  "
$ws.Range("B40").Value = "RenderToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.914"
$ws.Range("E40").Value = "  -5.81%  "
$ws.Range("B41").Value = "PaxDollar"
$ws.Range("C41").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.002"
$ws.Range("E41").Value = "  +0.65%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.519"
$ws.Range("E42").Value = "  -0.98%  "
$ws.Range("E43").Value = "  -3.91%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "7.004"
$ws.Range("E44").Value = "  -2.05%  "
$ws.Range("E45").Value = "  -2.37%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "59.81"
$ws.Range("E46").Value = "  -3.10%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "8.556"
$ws.Range("E47").Value = "  -3.25%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "32.90"
$ws.Range("E48").Value = "  -1.49%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.05507"
$ws.Range("E49").Value = "  -2.37%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.349"
$ws.Range("E50").Value = "  -4.52%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.3625"
$ws.Range("E51").Value = "  -3.83%  "
